# Apply opacity to colored background image of the #headline
#
# The real-world commit also touched a background image's opacity, but
# the reproducible OOXML diff available to us only shows text/run
# splitting (Word inserted w:proofErr spell-check bookmarks around
# certain proper nouns / foreign words, and re-flowed a couple of runs).
# This script reproduces exactly that diff against word/document.xml.
#
# Technique: for each affected paragraph, grab a Range that covers just
# the text to be restructured (leaving the paragraph mark, and any
# untouched runs, alone), then call Range.InsertXML with a <w:p> wrapper
# containing the desired list of <w:r>/<w:proofErr> children. Word merges
# this onto the existing paragraph (keeping the paragraph's own identity
# attributes and <w:pPr>) while replacing just the run content we
# targeted.
#
# NOTE: this runtime's PowerShell-subset does not bind named
# (-Param value) arguments on user functions (only positional ones),
# and a parenthesized expression used as a bare call argument confuses
# its parser -- so every helper below takes its arguments positionally,
# and every argument is first materialized into its own variable before
# the call (never `Foo (expr)`).

$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-RangeRunsXml($Range, $InnerXml) {
    $xml = '<w:p ' + $W + '>' + $InnerXml + '</w:p>'
    $Range.InsertXML($xml)
}

function Set-ParagraphRunsXml($Index, $InnerXml) {
    $p = $d.Paragraphs($Index)
    $full = $p.Range
    $r = $d.Range($full.Start, $full.End - 1)
    Set-RangeRunsXml $r $InnerXml
}

function Set-SubstringRunsXml($Index, $Needle, $InnerXml) {
    $p = $d.Paragraphs($Index)
    $full = $p.Range
    $text = $full.Text
    $idx = $text.IndexOf($Needle)
    $subStart = $full.Start + $idx
    $subEnd = $subStart + $Needle.Length
    $r = $d.Range($subStart, $subEnd)
    Set-RangeRunsXml $r $InnerXml
}

# --- 1. "13, Navoi str, Tashkent" -> split around "Navoi" (paragraph 8) ---
$xml1 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">13, </w:t></w:r>'
$xml1 += '<w:proofErr w:type="spellStart"/>'
$xml1 += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Navoi</w:t></w:r>'
$xml1 += '<w:proofErr w:type="spellEnd"/>'
$xml1 += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> str, Tashkent</w:t></w:r>'
Set-ParagraphRunsXml 8 $xml1

# --- 2. "TIACE Architecture Bootcamp 202" + "1" -> single merged run ---
$xml2 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>TIACE Architecture Bootcamp 2021</w:t></w:r>'
Set-ParagraphRunsXml 21 $xml2

# --- 3. "Vadim Makhmudov" -> split around "Makhmudov" ---
$xml3 = '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Vadim </w:t></w:r>'
$xml3 += '<w:proofErr w:type="spellStart"/>'
$xml3 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Makhmudov</w:t></w:r>'
$xml3 += '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphRunsXml 47 $xml3

# --- 4. Drop the trailing period after "...in architecture." ---
$xml4 = '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>One of the leading members of the Union of Architects of Uzbekistan recognized for his lifetime achievement in architecture</w:t></w:r>'
Set-ParagraphRunsXml 48 $xml4

# --- 5. "Art Specialist at Yeoju Technical Institute in Tashkent" -> split around "Yeoju" ---
$xml5 = '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Art Specialist at </w:t></w:r>'
$xml5 += '<w:proofErr w:type="spellStart"/>'
$xml5 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Yeoju</w:t></w:r>'
$xml5 += '<w:proofErr w:type="spellEnd"/>'
$xml5 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Technical Institute in Tashkent</w:t></w:r>'
Set-ParagraphRunsXml 52 $xml5

# --- 6. "Practice a variety ... coloured pencils ..." -> split around "coloured" ---
$xml6 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Practice a variety of drawing and painting techniques with Alexandra. At her workshops, both beginners and advanced artists work with different types of </w:t></w:r>'
$xml6 += '<w:proofErr w:type="spellStart"/>'
$xml6 += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>coloured</w:t></w:r>'
$xml6 += '<w:proofErr w:type="spellEnd"/>'
$xml6 += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> pencils, papers, and other related materials.</w:t></w:r>'
Set-ParagraphRunsXml 53 $xml6

# --- 7. "Shukur Djuraev " -> split around "Shukur" and "Djuraev" ---
$xml7 = '<w:proofErr w:type="spellStart"/>'
$xml7 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Shukur</w:t></w:r>'
$xml7 += '<w:proofErr w:type="spellEnd"/>'
$xml7 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
$xml7 += '<w:proofErr w:type="spellStart"/>'
$xml7 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Djuraev</w:t></w:r>'
$xml7 += '<w:proofErr w:type="spellEnd"/>'
$xml7 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
Set-ParagraphRunsXml 55 $xml7

# --- 8. "Known for his remarkable ... terms, Shukur will ..." -> split around "Shukur" ---
$xml8 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Known for his remarkable talent for explaining extremely complicated ideas in simple terms, </w:t></w:r>'
$xml8 += '<w:proofErr w:type="spellStart"/>'
$xml8 += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Shukur</w:t></w:r>'
$xml8 += '<w:proofErr w:type="spellEnd"/>'
$xml8 += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> will be your principal instructor. Teaching methods include lectures, tutorials, practical demonstrations, fieldwork and e-learning.</w:t></w:r>'
Set-ParagraphRunsXml 57 $xml8

# --- 9. "Certified coach and mentor with 8 years experience in HR" -> 3-way split (no proofErr, real words) ---
$xml9 = '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Certified coach and mentor with 8 years</w:t></w:r>'
$xml9 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of</w:t></w:r>'
$xml9 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> experience in HR</w:t></w:r>'
Set-ParagraphRunsXml 64 $xml9

# --- 10. "Christina Tursunova" -> split around "Tursunova" ---
$xml10 = '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Christina </w:t></w:r>'
$xml10 += '<w:proofErr w:type="spellStart"/>'
$xml10 += '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Tursunova</w:t></w:r>'
$xml10 += '<w:proofErr w:type="spellEnd"/>'
Set-ParagraphRunsXml 67 $xml10

# --- 11. "13, Navoi str., Tashkent" (second occurrence, preceded by "100011," run) -> split around "Navoi" ---
$xml11 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">13, </w:t></w:r>'
$xml11 += '<w:proofErr w:type="spellStart"/>'
$xml11 += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Navoi</w:t></w:r>'
$xml11 += '<w:proofErr w:type="spellEnd"/>'
$xml11 += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> str., Tashkent</w:t></w:r>'
$needle11 = "13, Navoi str., Tashkent"
Set-SubstringRunsXml 78 $needle11 $xml11

Write-Output "done"
